$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark cells as Text format first so numeric-looking values are not
# auto-converted to numbers by Excel (the source data are text strings).
foreach ($addr in @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.363.24"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").Value = "1.840.12"
$ws.Range("E3").Value = "  +4.10%  "
$ws.Range("D4").Value = "1.023"
$ws.Range("E4").Value = "  +2.42%  "
$ws.Range("D5").Value = "318.12"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").Value = "1.021"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("D7").Value = "0.4348"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").Value = "0.3721"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("D9").Value = "0.07338"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").Value = "21.39"
$ws.Range("E11").Value = "  +5.29%  "
$ws.Range("D12").Value = "2.029.81"
$ws.Range("E12").Value = "  +16.23%  "
$ws.Range("D13").Value = "5.476"
$ws.Range("E13").Value = "  +4.72%  "
$ws.Range("D14").Value = "6.671"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "0.07118"
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "82.05"
$ws.Range("E16").Value = "  +4.24%  "
$ws.Range("D17").Value = "1.027"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "0.000008999"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "15.42"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "27.404.43"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").Value = "11.13"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "2.207.23"
$ws.Range("E24").Value = "  +12.04%  "
$ws.Range("D25").Value = "156.68"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "1.905"
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("D27").Value = "18.55"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").Value = "5.300"
$ws.Range("E28").Value = "  +4.55%  "
$ws.Range("D29").Value = "1.927"
$ws.Range("E29").Value = "  +7.01%  "
$ws.Range("D30").Value = "115.48"
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "0.09018"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").Value = "1.202"
$ws.Range("E32").Value = "  +7.22%  "
$ws.Range("D33").Value = "0.7623"
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("D34").Value = "4.464"
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").Value = "2.834"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").Value = "1.023"
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").Value = "1.141"
$ws.Range("E37").Value = "  +3.53%  "
$ws.Range("D38").Value = "0.01954"
$ws.Range("E38").Value = "  +3.89%  "
$ws.Range("D39").Value = "0.05261"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5157"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.797"
$ws.Range("E41").Value = "  +7.57%  "
$ws.Range("D42").Value = "0.1664"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("D43").Value = "6.533"
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("D44").Value = "8.472"
$ws.Range("E44").Value = "  +5.96%  "
$ws.Range("D45").Value = "108.17"
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").Value = "10.52"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "1.023"
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("D48").Value = "0.4638"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").Value = "1.670"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").Value = "1.893"
$ws.Range("E50").Value = "  +9.42%  "
$ws.Range("D51").Value = "0.06279"
$ws.Range("E51").Value = "  +1.43%  "
